# Update Min Area Schedule
# - Uppercase the row-label strings in column A (Bedroom -> BEDROOM, etc.)
# - Move the active cell selection from L18 to E18
# - Set the page to portrait orientation (adds a <pageSetup> element)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column A labels, uppercased (order matters: it drives shared-string order)
$ws.Range("A2").Value = "BEDROOM"
$ws.Range("A3").Value = "BEDROOM 1"
$ws.Range("A4").Value = "BEDROOM 2"
$ws.Range("A5").Value = "BEDROOM 3"
$ws.Range("A6").Value = "BEDROOM 4"
$ws.Range("A7").Value = "STORAGE"
$ws.Range("A8").Value = "LIVING / DINING / KITCHEN"

# Move the selection / active cell to E18
$ws.Range("E18").Select()

# Page setup: portrait orientation
$ws.PageSetup.Orientation = 1
